# Update weekly Fruta/Hortaliza data rows (Papaya - Vega Central Mapocho de Santiago)
# The underlying data set was re-sorted/updated to reflect a new week's figures.
# Only columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de comercializacion),
# S (Precio $/Kg) and T (Kg / unidad) change for rows 2-11; every other column
# (Mercado, Region, Codreg, Tipo, Producto, Categoria, Variedad, Origen) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @{
    2  = @{ D = 44391; L = "Primera"; M = 15; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    3  = @{ D = 44391; L = "Segunda"; M = 20; N = 1000;  O = 1000;  P = 1000;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1000; T = 1 }
    4  = @{ D = 44343; L = "Primera"; M = 20; N = 1700;  O = 1700;  P = 1700;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1700; T = 1 }
    5  = @{ D = 44309; L = "Primera"; M = 10; N = 1600;  O = 1600;  P = 1600;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1600; T = 1 }
    6  = @{ D = 44195; L = "Primera"; M = 20; N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 10 kilos";          S = 1500; T = 10 }
    7  = @{ D = 44336; L = "Primera"; M = 10; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    8  = @{ D = 44292; L = "Primera"; M = 50; N = 14000; O = 14000; P = 14000; Q = "`$/bandeja 10 kilos";          S = 1400; T = 10 }
    9  = @{ D = 44400; L = "Primera"; M = 25; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    10 = @{ D = 44371; L = "Primera"; M = 20; N = 1800;  O = 1800;  P = 1800;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1800; T = 1 }
    11 = @{ D = 44371; L = "Segunda"; M = 30; N = 1200;  O = 1200;  P = 1200;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1200; T = 1 }
}

foreach ($row in $targets.Keys) {
    $vals = $targets[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
